# Apply the edits described by the diff via Find/Replace against the
# Word object model (Range.Find.Execute).
#
# Helper: replace the FIRST remaining match of $old with $new inside the
# given $range (ReplaceOne). Using MatchWholeWord:=$false avoids an
# observed quirk where whole-word matching on punctuation-only search
# strings (e.g. a run of dots) can select the wrong occurrence.
function Replace-First($range, [string]$old, [string]$new) {
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 0, $false, $new, 1) | Out-Null
}

$d = $word.ActiveDocument

# --- Paragraph 2: signature line (name + date) ------------------------
$p2 = $d.Paragraphs.Item(2)
Replace-First $p2.Range " Byś" " strózik"
$p2 = $d.Paragraphs.Item(2)
Replace-First $p2.Range "2021-03-03" "2021-03-11"

# --- Paragraph 5: address line -----------------------------------------
$p5 = $d.Paragraphs.Item(5)
Replace-First $p5.Range "ul. Lelewela 7/23, 01-476 Warszawa" "Matejki 18/6"

# --- Paragraph 20: "od dnia ... do ..." + destination town -------------
$p20 = $d.Paragraphs.Item(20)
Replace-First $p20.Range "2021-03-04" "2021-03-12"
$p20 = $d.Paragraphs.Item(20)
Replace-First $p20.Range "2021-03-04" "2021-03-14"
$p20 = $d.Paragraphs.Item(20)
Replace-First $p20.Range "Dębica" "Kraków"

# --- Paragraph 32: order number / order date placeholders --------------
$p32 = $d.Paragraphs.Item(32)
Replace-First $p32.Range " ............." " 20"
$p32 = $d.Paragraphs.Item(32)
Replace-First $p32.Range "............." "2021-03-11"

# --- Paragraph 33: "w dniach od ... do ..." -----------------------------
$p33 = $d.Paragraphs.Item(33)
Replace-First $p33.Range "2021-03-04" "2021-03-12"
$p33 = $d.Paragraphs.Item(33)
Replace-First $p33.Range "2021-03-04." "2021-03-14."

# --- Paragraph 56: statement date, destination, amount -----------------
$p56 = $d.Paragraphs.Item(56)
Replace-First $p56.Range "2021-03-03" "2021-03-11"
$p56 = $d.Paragraphs.Item(56)
Replace-First $p56.Range "Dębica" "Kraków"
$p56 = $d.Paragraphs.Item(56)
Replace-First $p56.Range "43.0 " "0.01 "
$p56 = $d.Paragraphs.Item(56)
Replace-First $p56.Range "czterdzieści trzy złote zero groszy" "zero złotych jeden grosz"
